$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.85133171081543
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = 2.579610109329224
$ws.Range("D1").Value = 1.088018536567688
$ws.Range("E1").Value = 0.741705596446991
